$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '43.452.88'
$ws.Range('E2').Value = '  +3.37%  '

# Row 3
$ws.Range('D3').Value = '2.309.18'
$ws.Range('E3').Value = '  +2.72%  '

# Row 4
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '105.49'
$ws.Range('E5').Value = '  +9.10%  '

# Row 6
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '308.42'
$ws.Range('E6').Value = '  +0.56%  '

# Row 7
$ws.Range('E7').Value = '  +0.48%  '

# Row 8
$ws.Range('E8').Value = '  -0.03%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.517'
$ws.Range('E9').Value = '  +5.68%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.04'
$ws.Range('E10').Value = '  +3.39%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.81'
$ws.Range('E11').Value = '  +2.32%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0812'
$ws.Range('E12').Value = '  -0.38%  '

# Row 13
$ws.Range('E13').Value = '  -1.43%  '

# Row 14
$ws.Range('E14').Value = '  +2.85%  '

# Row 15
$ws.Range('D15').Value = '2.665.19'
$ws.Range('E15').Value = '  +2.71%  '

# Row 16
$ws.Range('E16').Value = '  +4.70%  '

# Row 17
$ws.Range('D17').Value = '2.309.33'
$ws.Range('E17').Value = '  +2.97%  '

# Row 18
$ws.Range('E18').Value = '  +2.41%  '

# Row 19
$ws.Range('D19').Value = '43.375.32'
$ws.Range('E19').Value = '  +3.47%  '

# Row 20
$ws.Range('B20').Value = 'InternetComputer(DFINITY)'
$ws.Range('C20').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.90'
$ws.Range('E20').Value = '  -2.43%  '

# Row 21
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0₃0922'
$ws.Range('E21').Value = '  +2.31%  '

# Row 22
$ws.Range('E22').Value = '  +4.43%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.98'
$ws.Range('E23').Value = '  +1.33%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '240.91'
$ws.Range('E24').Value = '  +2.14%  '

# Row 25
$ws.Range('E25').Value = '  +4.47%  '

# Row 26
$ws.Range('E26').Value = '  +1.39%  '

# Row 27
$ws.Range('E27').Value = '  -0.58%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '25.05'
$ws.Range('E28').Value = '  +7.59%  '

# Row 29
$ws.Range('E29').Value = '  +4.69%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.33'
$ws.Range('E30').Value = '  -3.73%  '

# Row 31
$ws.Range('E31').Value = '  +1.23%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '162.29'
$ws.Range('E32').Value = '  -3.23%  '

# Row 33
$ws.Range('E33').Value = '  +1.26%  '

# Row 34
$ws.Range('E34').Value = '  +0.06%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.24'
$ws.Range('E35').Value = '  +4.48%  '

# Row 36
$ws.Range('E36').Value = '  +6.59%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0733'
$ws.Range('E37').Value = '  +1.92%  '

# Row 38
$ws.Range('E38').Value = '  +13.53%  '

# Row 39
$ws.Range('E39').Value = '  -1.01%  '

# Row 40
$ws.Range('E40').Value = '  +4.41%  '

# Row 41
$ws.Range('E41').Value = '  +2.49%  '

# Row 42
$ws.Range('E42').Value = '  +0.29%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.46'
$ws.Range('E43').Value = '  +13.32%  '

# Row 44
$ws.Range('E44').Value = '  +2.97%  '

# Row 45
$ws.Range('D45').Value = '1.961.83'
$ws.Range('E45').Value = '  +1.21%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '18.86'
$ws.Range('E46').Value = '  +2.12%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.04'
$ws.Range('E47').Value = '  +5.16%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.27'
$ws.Range('E48').Value = '  +6.58%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '58.10'
$ws.Range('E49').Value = '  +7.66%  '

# Row 50
$ws.Range('E50').Value = '  +1.57%  '

# Row 51
$ws.Range('E51').Value = '  +7.46%  '
